# correção nos dados e inicio da analise PNAD 2009
#
# The source table had "category header" rows (e.g. "sexo", "cor ou raça",
# "grupos de idade", "nível de instrução", "classes de rendimento mensal
# domiciliar per capita", "sem rendimento a menos ", the source footnote and
# the final footnote) that carried no numeric data. This edit removes those
# rows entirely (shifting the data rows below them up), which also causes
# the now-unused shared strings to be dropped when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row numbers (1-based, in the ORIGINAL layout) that only held a category
# label and no B:E data. Delete from bottom to top so row numbers of the
# rows still to be removed are not affected by earlier deletions.
$rowsToDelete = @(36, 35, 29, 27, 19, 13, 8, 5)

foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete() | Out-Null
}
